$d = $word.ActiveDocument

# Locate the paragraph that currently holds only the _GoBack bookmark
# (it is the second-to-last paragraph in the document body).
$count = $d.Paragraphs.Count
$bookmarkParaIndex = $count - 1

$bookmarkPara = $d.Paragraphs.Item($bookmarkParaIndex)
$r = $bookmarkPara.Range
$r.Collapse(1)
$r.InsertBefore("The problem is that the little girl is counting to ten on one hand and is probably confusing herself in the process. We are needing to find out where she will land if she continues to count this way on her left hand, the problem itself is unclear because I would think by reading it that if she counts this way every time then she would always end up on her first finger, but if she continues to count after reaching 10 on her first finger and moving to the middle finger for 11 then we have an issue, so that is what I will solve for. Our overall goal is to come up with an equation that will predict where the little girl will end up on 10, 100 and 1000.")

# Remove the final, otherwise-empty trailing paragraph by deleting the
# paragraph mark that used to separate it from the bookmark paragraph.
# This merges the (now empty) last paragraph away, leaving the bookmark
# paragraph (with its new text) as the final paragraph of the body.
$bookmarkPara2 = $d.Paragraphs.Item($bookmarkParaIndex)
$markRange = $d.Range($bookmarkPara2.Range.End - 1, $bookmarkPara2.Range.End)
$markRange.Delete()

# Apply numbered-list formatting ("List Paragraph" style + a new numbered
# list, matching the "1) 2) 3)" style already used elsewhere in the doc)
# to the paragraph that now holds the new text.
$finalPara = $d.Paragraphs.Item($bookmarkParaIndex)
$finalPara.Style = "List Paragraph"
$finalPara.Range.ListFormat.ApplyNumberDefault()

$lf = $finalPara.Range.ListFormat
$tmpl = $lf.ListTemplate

$lvl0 = $tmpl.ListLevels.Item(1)
$lvl0.NumberFormat = "%1)"

$lvl1 = $tmpl.ListLevels.Item(2)
$lvl1.NumberStyle = 4
$lvl1.NumberFormat = "%2."

$lvl2 = $tmpl.ListLevels.Item(3)
$lvl2.NumberStyle = 2
$lvl2.NumberFormat = "%3."

$lvl3 = $tmpl.ListLevels.Item(4)
$lvl3.NumberStyle = 0
$lvl3.NumberFormat = "%4."

$lvl4 = $tmpl.ListLevels.Item(5)
$lvl4.NumberStyle = 4
$lvl4.NumberFormat = "%5."

$lvl5 = $tmpl.ListLevels.Item(6)
$lvl5.NumberStyle = 2
$lvl5.NumberFormat = "%6."

$lvl6 = $tmpl.ListLevels.Item(7)
$lvl6.NumberStyle = 0
$lvl6.NumberFormat = "%7."

$lvl7 = $tmpl.ListLevels.Item(8)
$lvl7.NumberStyle = 4
$lvl7.NumberFormat = "%8."

$lvl8 = $tmpl.ListLevels.Item(9)
$lvl8.NumberStyle = 2
$lvl8.NumberFormat = "%9."
